$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 37037900
$ws.Range("I33").Value = 769.2353000000001
$ws.Range("J33").Value = 100001016
$ws.Range("K33").Value = 769.2353000000001
$ws.Range("L33").Value = 100001016
$ws.Range("M33").Value = -540.2353000000001
$ws.Range("N33").Value = -100001474
$ws.Range("H43").Value = 996.619
$ws.Range("J43").Value = 966.2632
$ws.Range("L43").Value = 966.2632
$ws.Range("N43").Value = -1104.2632
$ws.Range("H70").Value = 1542.0358
$ws.Range("J70").Value = 1693.7142
$ws.Range("L70").Value = 5081.142599999999
$ws.Range("N70").Value = -5621.142599999999
$ws.Range("H73").Value = 1542.0358
$ws.Range("J73").Value = 1693.7142
$ws.Range("L73").Value = 5081.142599999999
$ws.Range("N73").Value = -6953.142599999999
$ws.Range("H76").Value = 3043.4783
$ws.Range("J76").Value = 3266.6667
$ws.Range("L76").Value = 3266.6667
$ws.Range("N76").Value = -3896.6667
$ws.Range("H79").Value = 3043.4783
$ws.Range("J79").Value = 3266.6667
$ws.Range("L79").Value = 3266.6667
$ws.Range("N79").Value = -5450.6667
$ws.Range("H132").Value = 2282.2778
$ws.Range("I132").Value = 1227.119
$ws.Range("J132").Value = 5975.3335
$ws.Range("K132").Value = 3681.357
$ws.Range("L132").Value = 17926.0005
$ws.Range("M132").Value = -1151.357
$ws.Range("N132").Value = -22986.0005
$ws.Range("H135").Value = 905.4375
$ws.Range("I135").Value = 912.25
$ws.Range("J135").Value = 885
$ws.Range("K135").Value = 8210.25
$ws.Range("L135").Value = 7965
$ws.Range("M135").Value = -5675.25
$ws.Range("N135").Value = -13035
$ws.Range("H137").Value = 2169.1875
$ws.Range("I137").Value = 2026.3043
$ws.Range("J137").Value = 2534.3333
$ws.Range("K137").Value = 6078.9129
$ws.Range("L137").Value = 7602.999899999999
$ws.Range("M137").Value = -3528.9129
$ws.Range("N137").Value = -12702.9999
$ws.Range("H138").Value = 2600001
$ws.Range("I138").Value = 1195.3243
$ws.Range("J138").Value = 5003896.5
$ws.Range("K138").Value = 3585.9729
$ws.Range("L138").Value = 15011689.5
$ws.Range("M138").Value = 1554.0271
$ws.Range("N138").Value = -15021969.5
$ws.Range("H141").Value = 1497.95
$ws.Range("I141").Value = 1497.95
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4493.85
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 686.1499999999996
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 242.25
$ws.Range("I4").Value = 242.25
$ws.Range("K4").Value = 242.25
$ws.Range("M4").Value = -126.25
$ws.Range("H5").Value = 133
$ws.Range("I5").Value = 116.1
$ws.Range("K5").Value = 116.1
$ws.Range("M5").Value = -4.099999999999994
$ws.Range("H63").Value = 3274.125
$ws.Range("I63").Value = 2878.6
$ws.Range("J63").Value = 3933.3333
$ws.Range("K63").Value = 2878.6
$ws.Range("L63").Value = 3933.3333
$ws.Range("M63").Value = -2192.6
$ws.Range("N63").Value = -5305.3333
$ws.Range("H66").Value = 3274.125
$ws.Range("I66").Value = 2878.6
$ws.Range("J66").Value = 3933.3333
$ws.Range("K66").Value = 14393
$ws.Range("L66").Value = 19666.6665
$ws.Range("M66").Value = -10961
$ws.Range("N66").Value = -26530.6665
$ws.Range("H74").Value = 44792.695
$ws.Range("I74").Value = 77816.766
$ws.Range("J74").Value = 1861.4
$ws.Range("K74").Value = 77816.766
$ws.Range("L74").Value = 1861.4
$ws.Range("M74").Value = -76942.766
$ws.Range("N74").Value = -3609.4
$ws.Range("H77").Value = 44792.695
$ws.Range("I77").Value = 77816.766
$ws.Range("J77").Value = 1861.4
$ws.Range("K77").Value = 389083.83
$ws.Range("L77").Value = 9307
$ws.Range("M77").Value = -384715.83
$ws.Range("N77").Value = -18043

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 133
$ws.Range("I4").Value = 116.1
$ws.Range("K4").Value = 116.1
$ws.Range("M4").Value = -1.099999999999994
$ws.Range("H62").Value = 29900
$ws.Range("J62").Value = 29900
$ws.Range("L62").Value = 29900
$ws.Range("N62").Value = -31272
$ws.Range("H65").Value = 29900
$ws.Range("J65").Value = 29900
$ws.Range("L65").Value = 89700
$ws.Range("N65").Value = -96564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1280.9678
$ws.Range("I58").Value = 1180.381
$ws.Range("J58").Value = 1492.2
$ws.Range("K58").Value = 1180.381
$ws.Range("L58").Value = 1492.2
$ws.Range("M58").Value = -977.3810000000001
$ws.Range("N58").Value = -1898.2
$ws.Range("H134").Value = 1517.1936
$ws.Range("I134").Value = 1590.2084
$ws.Range("J134").Value = 1266.8572
$ws.Range("K134").Value = 4770.6252
$ws.Range("L134").Value = 3800.5716
$ws.Range("M134").Value = -2235.6252
$ws.Range("N134").Value = -8870.571599999999
$ws.Range("H136").Value = 1280.9678
$ws.Range("I136").Value = 1180.381
$ws.Range("J136").Value = 1492.2
$ws.Range("K136").Value = 3541.143
$ws.Range("L136").Value = 4476.6
$ws.Range("M136").Value = -991.143
$ws.Range("N136").Value = -9576.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 20000
$ws.Range("J49").Value = 20000
$ws.Range("L49").Value = 20000
$ws.Range("N49").Value = -20368
$ws.Range("H126").Value = 1971.4286
$ws.Range("I126").Value = 1971.4286
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5914.2858
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3444.2858
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3181.4666
$ws.Range("I132").Value = 2976.3809
$ws.Range("J132").Value = 3660
$ws.Range("K132").Value = 8929.1427
$ws.Range("L132").Value = 10980
$ws.Range("M132").Value = -6399.1427
$ws.Range("N132").Value = -16040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2243.5
$ws.Range("I40").Value = 1960.4
$ws.Range("J40").Value = 2951.25
$ws.Range("K40").Value = 1960.4
$ws.Range("L40").Value = 2951.25
$ws.Range("M40").Value = -1824.4
$ws.Range("N40").Value = -3223.25
$ws.Range("H55").Value = 605.5
$ws.Range("I55").Value = 179.1
$ws.Range("J55").Value = 1031.9
$ws.Range("K55").Value = 179.1
$ws.Range("L55").Value = 1031.9
$ws.Range("M55").Value = -6.099999999999994
$ws.Range("N55").Value = -1377.9
